$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.559.49'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.747.46'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'324.31"
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = "'0.4579"
$ws.Range('E7').Value = '  +9.04%  '
$ws.Range('D8').Value = "'0.3552"
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('D9').Value = "'0.07458"
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').Value = "'42.05"
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').Value = "'1.089"
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = "'1.003"
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').Value = "'20.65"
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').Value = "'5.969"
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = "'7.076"
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('D16').Value = '1.750.03'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = "'92.09"
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').Value = "'0.00001060"
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').Value = "'0.06413"
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').Value = "'1.002"
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = "'16.71"
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('D22').Value = "'5.775"
$ws.Range('E22').Value = '  -2.31%  '
$ws.Range('D23').Value = '27.620.66'
$ws.Range('D24').Value = "'11.16"
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = "'2.109"
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').Value = "'163.35"
$ws.Range('E26').Value = '  +4.17%  '
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').Value = '1.950.11'
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'125.78"
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = "'2.056"
$ws.Range('E30').Value = '  -3.18%  '
$ws.Range('D31').Value = "'1.049"
$ws.Range('E31').Value = '  -6.13%  '
$ws.Range('D32').Value = "'0.09198"
$ws.Range('E32').Value = '  +4.74%  '
$ws.Range('D33').Value = "'3.664"
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').Value = "'5.508"
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = "'0.02284"
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = "'11.72"
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('D37').Value = "'0.06022"
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').Value = "'0.2081"
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').Value = "'4.955"
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').Value = "'0.6282"
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').Value = "'1.202"
$ws.Range('E41').Value = '  +2.65%  '
$ws.Range('D42').Value = "'1.379"
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').Value = "'7.731"
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').Value = "'13.18"
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').Value = "'3.719"
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').Value = "'0.5860"
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('D47').Value = "'122.33"
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = "'1.931"
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.06855"
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = "'1.127"
$ws.Range('E50').Value = '  -3.90%  '
$ws.Range('D51').Value = "'71.67"
$ws.Range('E51').Value = '  -2.29%  '

# Reset style on cells that were forced to text via quote-prefix,
# so no stray number-format/quotePrefix style gets attached.
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
